# Auto-generated edit script: updates cryptos.xlsx Price (D) and
# Volume(1h) (E) columns with refreshed market data, and fixes the
# row order for two swapped coin pairs (rows 13/14 and 44/45/46/47).
#
# The Price/Volume columns are stored as literal TEXT in the source
# workbook (not numbers) - e.g. "28.675.60" or "17.00" or "0.07140" -
# so every write below is prefixed with a leading apostrophe to force
# Excel to keep treating the cell as text (preserving formats like
# trailing zeros) instead of auto-converting to a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = "'28.675.60"
$ws.Range('E2').Value = "'  -2.73%  "

# Row 3
$ws.Range('D3').Value = "'1.852.06"
$ws.Range('E3').Value = "'  -3.47%  "

# Row 4
$ws.Range('E4').Value = "'  -0.68%  "

# Row 5
$ws.Range('D5').Value = "'335.61"
$ws.Range('E5').Value = "'  +3.07%  "

# Row 6
$ws.Range('E6').Value = "'  -0.66%  "

# Row 7
$ws.Range('D7').Value = "'0.4645"
$ws.Range('E7').Value = "'  -3.55%  "

# Row 8
$ws.Range('D8').Value = "'0.3927"
$ws.Range('E8').Value = "'  -3.13%  "

# Row 9
$ws.Range('D9').Value = "'46.57"
$ws.Range('E9').Value = "'  -2.70%  "

# Row 10
$ws.Range('D10').Value = "'0.07937"
$ws.Range('E10').Value = "'  -3.36%  "

# Row 11
$ws.Range('D11').Value = "'0.9836"
$ws.Range('E11').Value = "'  -2.37%  "

# Row 12
$ws.Range('E12').Value = "'  -4.77%  "

# Row 13
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = "'5.838"
$ws.Range('E13').Value = "'  -3.55%  "

# Row 14
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = "'1.826.55"
$ws.Range('E14').Value = "'  -4.99%  "

# Row 15
$ws.Range('D15').Value = "'7.017"
$ws.Range('E15').Value = "'  -2.87%  "

# Row 16
$ws.Range('D16').Value = "'0.06801"
$ws.Range('E16').Value = "'  -0.74%  "

# Row 17
$ws.Range('E17').Value = "'  -0.69%  "

# Row 18
$ws.Range('D18').Value = "'87.55"
$ws.Range('E18').Value = "'  -4.09%  "

# Row 19
$ws.Range('E19').Value = "'  -2.31%  "

# Row 20
$ws.Range('D20').Value = "'17.00"
$ws.Range('E20').Value = "'  -2.87%  "

# Row 21
$ws.Range('D21').Value = "'1.003"
$ws.Range('E21').Value = "'  -0.65%  "

# Row 22
$ws.Range('D22').Value = "'28.667.60"
$ws.Range('E22').Value = "'  -2.77%  "

# Row 23
$ws.Range('D23').Value = "'5.407"
$ws.Range('E23').Value = "'  -4.49%  "

# Row 24
$ws.Range('D24').Value = "'11.36"
$ws.Range('E24').Value = "'  -4.21%  "

# Row 25
$ws.Range('E25').Value = "'  -2.74%  "

# Row 26
$ws.Range('D26').Value = "'2.043.63"
$ws.Range('E26').Value = "'  -5.25%  "

# Row 27
$ws.Range('D27').Value = "'153.22"
$ws.Range('E27').Value = "'  -1.87%  "

# Row 28
$ws.Range('D28').Value = "'6.253"
$ws.Range('E28').Value = "'  -5.20%  "

# Row 29
$ws.Range('D29').Value = "'19.43"
$ws.Range('E29').Value = "'  -2.59%  "

# Row 30
$ws.Range('D30').Value = "'2.036"
$ws.Range('E30').Value = "'  -2.96%  "

# Row 31
$ws.Range('D31').Value = "'117.36"
$ws.Range('E31').Value = "'  -2.60%  "

# Row 32
$ws.Range('D32').Value = "'0.9846"
$ws.Range('E32').Value = "'  -2.69%  "

# Row 33
$ws.Range('D33').Value = "'0.09424"
$ws.Range('E33').Value = "'  -2.06%  "

# Row 34
$ws.Range('D34').Value = "'5.391"
$ws.Range('E34').Value = "'  -3.93%  "

# Row 35
$ws.Range('D35').Value = "'3.498"
$ws.Range('E35').Value = "'  -1.59%  "

# Row 36
$ws.Range('D36').Value = "'1.351"
$ws.Range('E36').Value = "'  -1.44%  "

# Row 37
$ws.Range('D37').Value = "'0.06152"
$ws.Range('E37').Value = "'  -2.21%  "

# Row 38
$ws.Range('D38').Value = "'0.02204"
$ws.Range('E38').Value = "'  -3.45%  "

# Row 39
$ws.Range('D39').Value = "'1.163"
$ws.Range('E39').Value = "'  -1.33%  "

# Row 40
$ws.Range('E40').Value = "'  -3.36%  "

# Row 41
$ws.Range('D41').Value = "'7.636"
$ws.Range('E41').Value = "'  -3.06%  "

# Row 42
$ws.Range('D42').Value = "'10.12"
$ws.Range('E42').Value = "'  -5.68%  "

# Row 43
$ws.Range('E43').Value = "'  -3.17%  "

# Row 44
$ws.Range('B44').Value = 'WEMIXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D44').Value = "'1.254"
$ws.Range('E44').Value = "'  -2.12%  "

# Row 45
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').Value = "'2.348"
$ws.Range('E45').Value = "'  -2.14%  "

# Row 46
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').Value = "'0.5415"
$ws.Range('E46').Value = "'  -2.64%  "

# Row 47
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = "'11.84"
$ws.Range('E47').Value = "'  -4.71%  "

# Row 48
$ws.Range('D48').Value = "'0.07140"
$ws.Range('E48').Value = "'  -4.39%  "

# Row 49
$ws.Range('D49').Value = "'1.918"
$ws.Range('E49').Value = "'  -0.76%  "

# Row 50
$ws.Range('D50').Value = "'115.79"
$ws.Range('E50').Value = "'  -1.96%  "

# Row 51
$ws.Range('D51').Value = "'43.56"
$ws.Range('E51').Value = "'  +3.91%  "
